$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "755÷7=107, 6" "468÷4=117, 0"
Replace-Text "268÷2=134, 0" "382÷2=191, 0"
Replace-Text "651÷8=81, 3" "753÷5=150, 3"
Replace-Text "340÷4=85, 0" "901÷8=112, 5"
Replace-Text "930÷4=232, 2" "816÷9=90, 6"
Replace-Text "530÷3=176, 2" "643÷3=214, 1"
Replace-Text "489÷7=69, 6" "160÷9=17, 7"
Replace-Text "849÷9=94, 3" "590÷3=196, 2"
Replace-Text "586÷3=195, 1" "540÷8=67, 4"
Replace-Text "732÷7=104, 4" "974÷4=243, 2"
Replace-Text "657÷5=131, 2" "125÷6=20, 5"
Replace-Text "630÷6=105, 0" "184÷8=23, 0"
Replace-Text "213÷3=71, 0" "748÷2=374, 0"
Replace-Text "881÷8=110, 1" "146÷2=73, 0"
Replace-Text "581÷9=64, 5" "386÷5=77, 1"
Replace-Text "407÷4=101, 3" "920÷9=102, 2"
Replace-Text "733÷6=122, 1" "675÷8=84, 3"
Replace-Text "211÷8=26, 3" "804÷5=160, 4"
Replace-Text "642÷4=160, 2" "673÷9=74, 7"
Replace-Text "463÷3=154, 1" "869÷4=217, 1"
Replace-Text "410÷6=68, 2" "458÷6=76, 2"
Replace-Text "823÷8=102, 7" "100÷5=20, 0"
Replace-Text "698÷6=116, 2" "971÷5=194, 1"
Replace-Text "418÷5=83, 3" "385÷6=64, 1"
Replace-Text "354÷8=44, 2" "229÷7=32, 5"
